$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, pushing existing rows 217-279 down to 218-280
$ws.Rows(217).Insert()

# Populate the newly inserted row 217 with the new data record
$ws.Cells(217,1).Value  = 9
$ws.Cells(217,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells(217,3).Value  = "Metropolitana"
$ws.Cells(217,4).Value  = 44559
$ws.Cells(217,5).Value  = 13
$ws.Cells(217,6).Value  = 100112032
$ws.Cells(217,7).Value  = "Zapallo italiano"
$ws.Cells(217,8).Value  = "Sin especificar"
$ws.Cells(217,9).Value  = "Primera"
$ws.Cells(217,10).Value = 160
$ws.Cells(217,11).Value = 7000
$ws.Cells(217,12).Value = 8000
$ws.Cells(217,13).Value = 7500
$ws.Cells(217,14).Value = "`$/caja 60 unidades"
$ws.Cells(217,15).Value = "Región Metropolitana"
$ws.Cells(217,16).Value = 125
$ws.Cells(217,17).Value = 60
$ws.Cells(217,18).Value = "Hortaliza"
